$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.938.89"
$ws.Range("E2").Value = "  -0.71%  "
$ws.Range("D3").Value = "2.525.63"
$ws.Range("E3").Value = "  +2.81%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "'540.01"
$ws.Range("E5").Value = "  +0.60%  "
$ws.Range("D6").Value = "'143.26"
$ws.Range("E6").Value = "  -3.28%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").Value = "'0.571"
$ws.Range("E8").Value = "  +0.54%  "
$ws.Range("D9").Value = "2.524.82"
$ws.Range("E9").Value = "  +2.07%  "
$ws.Range("E10").Value = "  +0.90%  "
$ws.Range("E11").Value = "  +0.31%  "
$ws.Range("E12").Value = "  +4.57%  "
$ws.Range("E13").Value = "  +0.60%  "
$ws.Range("D14").Value = "2.967.71"
$ws.Range("E14").Value = "  +2.14%  "
$ws.Range("D15").Value = "'23.49"
$ws.Range("E15").Value = "  -2.61%  "
$ws.Range("D16").Value = "58.928.52"
$ws.Range("E16").Value = "  -0.69%  "
$ws.Range("E17").Value = "  +1.06%  "
$ws.Range("D18").Value = "2.517.91"
$ws.Range("E18").Value = "  +1.00%  "
$ws.Range("D19").Value = "'11.21"
$ws.Range("E19").Value = "  +0.33%  "
$ws.Range("E20").Value = "  -1.77%  "
$ws.Range("D21").Value = "'324.66"
$ws.Range("E21").Value = "  +0.15%  "
$ws.Range("D22").Value = "'0.999"
$ws.Range("E22").Value = "  +3.59%  "
$ws.Range("E23").Value = "  +0.72%  "
$ws.Range("D24").Value = "'62.08"
$ws.Range("E24").Value = "  +2.37%  "
$ws.Range("D25").Value = "'0.440"
$ws.Range("E25").Value = "  -4.89%  "
$ws.Range("E26").Value = "  +0.78%  "
$ws.Range("D27").Value = "2.619.80"
$ws.Range("E27").Value = "  +2.00%  "
$ws.Range("D28").Value = "'0.995"
$ws.Range("E28").Value = "  +1.76%  "
$ws.Range("D29").Value = "'7.78"
$ws.Range("E29").Value = "  +0.59%  "
$ws.Range("D30").Value = "0.0₃0774"
$ws.Range("E30").Value = "  +0.11%  "
$ws.Range("D31").Value = "'1.81"
$ws.Range("E31").Value = "  -0.39%  "
$ws.Range("D32").Value = "'6.66"
$ws.Range("E32").Value = "  -2.34%  "
$ws.Range("D33").Value = "'1.19"
$ws.Range("E33").Value = "  -6.12%  "
$ws.Range("D34").Value = "'0.998"
$ws.Range("E34").Value = "  +0.02%  "
$ws.Range("B35").Value = "Monero"
$ws.Range("C35").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D35").Value = "'156.27"
$ws.Range("E35").Value = "  +0.49%  "
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").Value = "'1.44"
$ws.Range("E36").Value = "  +3.41%  "
$ws.Range("D37").Value = "'18.64"
$ws.Range("E37").Value = "  +1.26%  "
$ws.Range("E38").Value = "  -4.89%  "
$ws.Range("E39").Value = "  -9.49%  "
$ws.Range("D40").Value = "'5.69"
$ws.Range("E40").Value = "  -4.06%  "
$ws.Range("D41").Value = "'36.96"
$ws.Range("E41").Value = "  +0.57%  "
$ws.Range("D42").Value = "'295.63"
$ws.Range("E42").Value = "  -7.26%  "
$ws.Range("E43").Value = "  -0.62%  "
$ws.Range("D44").Value = "'0.819"
$ws.Range("E44").Value = "  -2.92%  "
$ws.Range("D45").Value = "'0.998"
$ws.Range("E45").Value = "  +0.24%  "
$ws.Range("D46").Value = "'0.598"
$ws.Range("E46").Value = "  +2.25%  "
$ws.Range("D47").Value = "'10.80"
$ws.Range("E47").Value = "  +0.54%  "
$ws.Range("D48").Value = "'0.0929"
$ws.Range("E48").Value = "  -1.15%  "
$ws.Range("D49").Value = "'122.62"
$ws.Range("E49").Value = "  +0.80%  "
$ws.Range("D50").Value = "'18.55"
$ws.Range("E50").Value = "  +0.02%  "
$ws.Range("E51").Value = "  -0.32%  "
